# Reproduce the OPOST_AUTOMATION spreadsheet update:
#  - clear the old text-based style (numFmtId 49 "@") that was inherited from the
#    column style so the table body can go back to a clean "Normal"/General look
#  - give the table a bigger (20pt) font
#  - turn the "is_random" column (D) into a real numeric column (0.00 format)
#    instead of text ("1"/"1") so D2/D3 become actual numbers (1 and 0)
#  - add a 4th (blank) formatted row under the existing two data rows, and
#    pre-format D5:D17 with the same numeric format for future rows
#  - select D2, matching the author's last selection when the file was saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- formatting -----------------------------------------------------------

# Reset A1:D4 back to the plain "Normal" style first. This clears the text
# ("@") number format that columns C:D inherited from the column style, so
# we can rebuild the formatting cleanly (matches cellXfs entries 2-4).
$ws.Range("A1:D4").Style = "Normal"

# Whole table (A:C) gets a larger font.
$ws.Range("A1:C4").Font.Size = 20

# Column D (is_random) also gets the larger font for the header/data rows,
# plus a numeric "0.00" format so true numbers (not text) are stored.
$ws.Range("D1:D4").Font.Size = 20
$ws.Range("D1:D4").NumberFormat = "0.00"

# Pre-format the rest of column D (rows 5-17) with the same numeric format,
# using the regular (default-size) font.
$ws.Range("D5:D17").NumberFormat = "0.00"

# Give the header/data rows some extra height to match the bigger font.
$ws.Range("A1:D4").RowHeight = 25.8

# --- values -----------------------------------------------------------

# D2/D3 used to store the text "1" as a shared string; now store real
# numbers instead (1 for ديما, 0 for جيانا). NumberFormat was already
# switched away from Text above, so these are written as numeric cells.
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 0

# --- selection --------------------------------------------------------

[void]$ws.Range("D2").Select()
